$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row pairs that were swapped (entire row content swapped between the two rows).
$pairs = @(
    @(3, 4),
    @(6, 7),
    @(13, 14),
    @(16, 17),
    @(22, 23)
)

# Columns A (1) through AY (51). Skip Y (25) and AA (27) -- text-looking dates that
# are identical within every swapped pair, so skipping avoids Excel's automatic
# date-string-to-serial conversion on round-trip while producing the same result.
$skipCols = @(25, 27)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    for ($col = 1; $col -le 51; $col++) {
        if ($skipCols -contains $col) { continue }

        $cell1 = $ws.Cells.Item($r1, $col)
        $cell2 = $ws.Cells.Item($r2, $col)

        $v1 = $cell1.Value2
        $v2 = $cell2.Value2

        $cell1.Value2 = $v2
        $cell2.Value2 = $v1
    }
}
